# Updated cryptos list on Wed Nov  1 21:56:11 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to remain plain text even when the string parses as a
    # number (e.g. "228.05"), then drop the temporary number-format style
    # so the cell's style index is unchanged from its original (default).
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "35.286.97"
$ws.Range("E2").Value = "  +2.16%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.840.04"
$ws.Range("E3").Value = "  +1.81%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
Set-TextValue "D5" "228.05"
$ws.Range("E5").Value = "  +1.07%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +2.44%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.06%  "

# Row 8 - Solana
Set-TextValue "D8" "43.34"
$ws.Range("E8").Value = "  +15.89%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.306"
$ws.Range("E9").Value = "  +5.15%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.0694"
$ws.Range("E10").Value = "  +1.85%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +3.66%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.107.01"

# Row 13 - Chainlink
Set-TextValue "D13" "11.61"
$ws.Range("E13").Value = "  +2.88%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.839.48"
$ws.Range("E14").Value = "  +1.61%  "

# Row 15 - Polkadot
Set-TextValue "D15" "4.76"
$ws.Range("E15").Value = "  +7.58%  "

# Row 16 - Polygon
Set-TextValue "D16" "0.661"
$ws.Range("E16").Value = "  +4.55%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "35.201.90"
$ws.Range("E17").Value = "  +2.16%  "

# Row 18 - Litecoin
Set-TextValue "D18" "70.08"
$ws.Range("E18").Value = "  +2.19%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "245.95"
$ws.Range("E19").Value = "  +1.11%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0795"
$ws.Range("E20").Value = "  +2.62%  "

# Row 21 - Avalanche
Set-TextValue "D21" "12.16"
$ws.Range("E21").Value = "  +8.79%  "

# Row 22 - Uniswap
Set-TextValue "D22" "4.70"
$ws.Range("E22").Value = "  +13.91%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.12%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -0.47%  "

# Row 25 - Monero
Set-TextValue "D25" "172.43"
$ws.Range("E25").Value = "  +0.60%  "

# Row 26 - Cosmos
Set-TextValue "D26" "7.94"
$ws.Range("E26").Value = "  +1.58%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "17.92"
$ws.Range("E27").Value = "  +3.61%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  +2.06%  "

# Row 29 - EURNeutrino
$ws.Range("D29").Value = "3.614.99"
$ws.Range("E29").Value = "  +48.78%  "

# Row 31 - PancakeSwap
Set-TextValue "D31" "1.32"
$ws.Range("E31").Value = "  +7.35%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +3.48%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +3.75%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  +3.89%  "

# Row 35 - LidoDAOToken
Set-TextValue "D35" "1.89"
$ws.Range("E35").Value = "  +4.25%  "

# Row 36 - ImmutableX
Set-TextValue "D36" "0.676"
$ws.Range("E36").Value = "  +3.64%  "

# Row 37 - Aave
Set-TextValue "D37" "90.43"
$ws.Range("E37").Value = "  +12.17%  "

# Rows 38 and 39 swap places: TrustWalletToken <-> Maker (values also updated)
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.344.71"
$ws.Range("E38").Value = "  -1.34%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D39" "1.08"
$ws.Range("E39").Value = "  +1.17%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  +9.10%  "

# Row 41 - RenderToken
Set-TextValue "D41" "2.44"
$ws.Range("E41").Value = "  +3.31%  "

# Row 42 - VeChain
$ws.Range("E42").Value = "  +4.01%  "

# Row 43 - InjectiveProtocol
Set-TextValue "D43" "14.93"
$ws.Range("E43").Value = "  +8.82%  "

# Row 44 - WEMIXToken
$ws.Range("E44").Value = "  +5.10%  "

# Row 45 - HuobiToken
$ws.Range("E45").Value = "  +1.05%  "

# Row 46 - MXToken
Set-TextValue "D46" "2.82"
$ws.Range("E46").Value = "  +1.70%  "

# Row 47 - Kaspa
$ws.Range("E47").Value = "  +3.52%  "

# Row 48 - FraxShare
Set-TextValue "D48" "6.08"
$ws.Range("E48").Value = "  +4.83%  "

# Row 49 - RocketPoolETH
$ws.Range("D49").Value = "2.007.32"
$ws.Range("E49").Value = "  +1.95%  "

# Row 50 - Quant
Set-TextValue "D50" "104.74"
$ws.Range("E50").Value = "  +2.13%  "

# Row 51 - PaxDollar
$ws.Range("E51").Value = "  -0.13%  "
